# Scheduled-runner refresh of market-board derived Leve profit figures
# (currentAveragePrice / LevePrice* / LeveProfit* columns) across all
# eight crafting-job sheets. Values below were produced by an external
# price refresh; only numeric H:N cells change (no formulas involved).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1237.711
$ws.Range("I15").Value = 1237.711
$ws.Range("K15").Value = 3713.133
$ws.Range("M15").Value = -3544.133

$ws.Range("H19").Value = 398
$ws.Range("I19").Value = 132.71428
$ws.Range("K19").Value = 132.71428
$ws.Range("M19").Value = 42.28572

$ws.Range("H58").Value = 3411.1924
$ws.Range("J58").Value = 6806
$ws.Range("L58").Value = 20418
$ws.Range("N58").Value = -20718

$ws.Range("H70").Value = 551909.2
$ws.Range("I70").Value = 835266.25
$ws.Range("K70").Value = 2505798.75
$ws.Range("M70").Value = -2505528.75

$ws.Range("H73").Value = 551909.2
$ws.Range("I73").Value = 835266.25
$ws.Range("K73").Value = 2505798.75
$ws.Range("M73").Value = -2504862.75

$ws.Range("H137").Value = 2742.5254
$ws.Range("I137").Value = 2213.1064
$ws.Range("K137").Value = 6639.3192
$ws.Range("M137").Value = -4089.3192

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3188.281
$ws.Range("I32").Value = 3188.281
$ws.Range("K32").Value = 3188.281
$ws.Range("M32").Value = -2901.281

$ws.Range("H45").Value = 2639.7856
$ws.Range("I45").Value = 2041.25
$ws.Range("J45").Value = 3437.8333
$ws.Range("K45").Value = 2041.25
$ws.Range("L45").Value = 3437.8333
$ws.Range("M45").Value = -1664.25
$ws.Range("N45").Value = -4191.8333

$ws.Range("H74").Value = 1738.5079
$ws.Range("I74").Value = 1754.8704
$ws.Range("K74").Value = 1754.8704
$ws.Range("M74").Value = -880.8704

$ws.Range("H77").Value = 1738.5079
$ws.Range("I77").Value = 1754.8704
$ws.Range("K77").Value = 8774.352000000001
$ws.Range("M77").Value = -4406.352000000001

$ws.Range("H122").Value = 5479.933
$ws.Range("I122").Value = 1949.5
$ws.Range("J122").Value = 6023.077
$ws.Range("K122").Value = 5848.5
$ws.Range("L122").Value = 18069.231
$ws.Range("M122").Value = -3398.5
$ws.Range("N122").Value = -22969.231

$ws.Range("H132").Value = 3102.314
$ws.Range("I132").Value = 2542.2754
$ws.Range("J132").Value = 5375.4116
$ws.Range("K132").Value = 7626.8262
$ws.Range("L132").Value = 16126.2348
$ws.Range("M132").Value = -5096.8262
$ws.Range("N132").Value = -21186.2348

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 35685.715
$ws.Range("J60").Value = 39966.668
$ws.Range("L60").Value = 39966.668
$ws.Range("N60").Value = -40988.668

$ws.Range("H107").Value = 491
$ws.Range("I107").Value = 291.85715
$ws.Range("K107").Value = 291.85715
$ws.Range("M107").Value = 1628.14285

$ws.Range("H119").Value = 75000
$ws.Range("J119").Value = 75000
$ws.Range("L119").Value = 75000
$ws.Range("N119").Value = -84676

$ws.Range("H132").Value = 3560.3704
$ws.Range("I132").Value = 2038.421
$ws.Range("K132").Value = 6115.263
$ws.Range("M132").Value = -3585.263

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 34666860
$ws.Range("I11").Value = 52000160
$ws.Range("J11").Value = 259
$ws.Range("K11").Value = 156000480
$ws.Range("L11").Value = 777
$ws.Range("M11").Value = -156000340
$ws.Range("N11").Value = -1057

$ws.Range("H80").Value = 3312.3333
$ws.Range("J80").Value = 2972.5
$ws.Range("L80").Value = 8917.5
$ws.Range("N80").Value = -10789.5

$ws.Range("H83").Value = 3312.3333
$ws.Range("J83").Value = 2972.5
$ws.Range("L83").Value = 26752.5
$ws.Range("N83").Value = -36112.5

$ws.Range("H92").Value = 1187.3334
$ws.Range("J92").Value = 1198.5385
$ws.Range("L92").Value = 3595.6155
$ws.Range("N92").Value = -6091.6155

$ws.Range("H93").Value = 6692.6665
$ws.Range("I93").Value = 6080
$ws.Range("K93").Value = 18240
$ws.Range("M93").Value = -16368

$ws.Range("H98").Value = 3400.2104
$ws.Range("I98").Value = 3025.9092
$ws.Range("J98").Value = 3914.875
$ws.Range("K98").Value = 9077.7276
$ws.Range("L98").Value = 11744.625
$ws.Range("M98").Value = -7579.7276
$ws.Range("N98").Value = -14740.625

$ws.Range("H110").Value = 3030
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 3030
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 9090
$ws.Range("M110").Value = ""
$ws.Range("N110").Value = -17270

$ws.Range("H122").Value = 72699.57000000001
$ws.Range("J122").Value = 144423.14
$ws.Range("L122").Value = 1299808.26
$ws.Range("N122").Value = -1304708.26

$ws.Range("H129").Value = 68453.92999999999
$ws.Range("I129").Value = 773.6
$ws.Range("J129").Value = 102294.1
$ws.Range("K129").Value = 2320.8
$ws.Range("L129").Value = 306882.3
$ws.Range("M129").Value = 2679.2
$ws.Range("N129").Value = -316882.3

$ws.Range("H139").Value = 6557.7354
$ws.Range("I139").Value = 3369.0715
$ws.Range("K139").Value = 10107.2145
$ws.Range("M139").Value = -4967.2145

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H130").Value = 69997
$ws.Range("J130").Value = 69997
$ws.Range("L130").Value = 69997
$ws.Range("N130").Value = -80037

$ws.Range("H132").Value = 46074.754
$ws.Range("I132").Value = 6087.6855
$ws.Range("J132").Value = 146042.42
$ws.Range("K132").Value = 18263.0565
$ws.Range("L132").Value = 438127.26
$ws.Range("M132").Value = -15733.0565
$ws.Range("N132").Value = -443187.26

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7579.4
$ws.Range("I7").Value = 6000
$ws.Range("K7").Value = 6000
$ws.Range("M7").Value = -5888

$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").Value = ""

$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").Value = ""

$ws.Range("H94").Value = 40993.5
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 40993.5
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 40993.5
$ws.Range("M94").Value = ""
$ws.Range("N94").Value = -42345.5

$ws.Range("H126").Value = 7579.4
$ws.Range("I126").Value = 6000
$ws.Range("K126").Value = 18000
$ws.Range("M126").Value = -15530

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 31000
$ws.Range("J54").Value = 31000
$ws.Range("L54").Value = 31000
$ws.Range("N54").Value = -32040

$ws.Range("H113").Value = 1065.3572
$ws.Range("I113").Value = 948.36365
$ws.Range("K113").Value = 2845.09095
$ws.Range("M113").Value = -675.0909499999998

$ws.Range("H132").Value = 34698.535
$ws.Range("I132").Value = 1454.5217
$ws.Range("J132").Value = 143928.86
$ws.Range("K132").Value = 4363.5651
$ws.Range("L132").Value = 431786.58
$ws.Range("M132").Value = -1833.5651
$ws.Range("N132").Value = -436846.58

$ws.Range("H136").Value = 275823.47
$ws.Range("I136").Value = 280981.06
$ws.Range("K136").Value = 842943.1799999999
$ws.Range("M136").Value = -840393.1799999999
